# Add a basic BOM list (rows 10-30) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "PPG Module connector"
$ws.Range("C10").Value = "U`$1"
$ws.Range("D10").Value = "Digi-key Part number: S9014E-50-ND"

# Row 11
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "10k resistor"
$ws.Range("C11").Value = "R12, R10, R11"

# Row 12
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "TIA"
$ws.Range("C12").Value = "TIAG`$1"

# Row 13
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "1uF capacitor"
$ws.Range("C13").Value = "C14"

# Row 14
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "unknown resistor"
$ws.Range("C14").Value = "R9"

# Row 15
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "crystal"
$ws.Range("C15").Value = "Y1"

# Row 16
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "100nF capacitor"
$ws.Range("C16").Value = "C3, C6, C4, C10, C13, C12, C11"

# Row 17
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "2.2uF capacitor"
$ws.Range("C17").Value = "C8, C9"

# Row 18
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "CONN_07-1.27MM"
$ws.Range("C18").Value = "J1"

# Row 19
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "220k resistor"
$ws.Range("C19").Value = "R1, R3, R5, R6"

# Row 20
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "1M resistor"
$ws.Range("C20").Value = "R2, R8"

# Row 21
$ws.Range("A21").Value = 12
$ws.Range("B21").Value = "10M resistor"
$ws.Range("C21").Value = "R4, R7"

# Row 22
$ws.Range("A22").Value = 13
$ws.Range("B22").Value = "2.2nF capactior"
$ws.Range("C22").Value = "C7, C1"

# Row 23
$ws.Range("A23").Value = 14
$ws.Range("B23").Value = "Op Amp"
$ws.Range("C23").Value = "OP2G`$2, OP2G`$3, OP2G`$1, OP2G`$4"

# Row 24
$ws.Range("A24").Value = 15
$ws.Range("B24").Value = "22uF capacitor"
$ws.Range("C24").Value = "C2"

# Row 25
$ws.Range("A25").Value = 16
$ws.Range("B25").Value = "22nF capacitor"
$ws.Range("C25").Value = "C5"

# Row 26
$ws.Range("A26").Value = 17
$ws.Range("B26").Value = "AD5242"
$ws.Range("C26").Value = "AD5242"

# Row 27
$ws.Range("A27").Value = 18
$ws.Range("B27").Value = "AD5171"
$ws.Range("C27").Value = "AD5171"

# Row 28
$ws.Range("A28").Value = 19
$ws.Range("B28").Value = "Trasnsitor"
$ws.Range("C28").Value = "Q2"
$ws.Range("D28").Value = "NPN-GENERIC"

# Row 29
$ws.Range("A29").Value = 20
$ws.Range("B29").Value = "LSM6DS3"
$ws.Range("C29").Value = "IMU1"
$ws.Range("D29").Value = "LSM6DS3_CAST"

# Row 30
$ws.Range("A30").Value = 21
$ws.Range("B30").Value = "PCF8523"
$ws.Range("C30").Value = "U1"

# Update the active cell selection to match the author's final cursor position.
$null = $ws.Range("E20").Select()
